# This script re-applies a cyclic re-shuffle of the species-observation
# records that live in rows 22-25 and a swap between rows 29-30 of the
# "Artfynd" sheet. Each physical spreadsheet row keeps its location-derived
# columns (P, S, T, U, V, W, Y, Z, AA, AB, AD, AE, AG, AT, AW, AY, ...) but
# the species identity columns (A, B, E, F, G, H), the coordinates (Q, R),
# the observer order (AX) and the handful of "species-shape" blank cells
# (J vs L/M, and AF) move with the record.
#
# Net effect (old row -> new row content):
#   row22 <- old row23 data
#   row23 <- old row24 data
#   row24 <- old row25 data
#   row25 <- old row22 data
#   row29 <- old row30 data
#   row30 <- old row29 data

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Txt($addr, $text) {
    $ws.Range($addr).Value = $text
}

function Clear-Cell($addr) {
    $ws.Range($addr).ClearContents()
}

# ---- Row 22 (becomes the old row23 record) ----
Set-Txt "A22" 111941777
Set-Txt "B22" 77515
Set-Txt "E22" 6425
Set-Txt "F22" "Garnlav"
Set-Txt "G22" "Alectoria sarmentosa"
Set-Txt "H22" "(Ach.) Ach."
Clear-Cell "L22"
Clear-Cell "M22"
Set-Txt "J22" ""
Set-Txt "AF22" ""
Set-Txt "Q22" 466414.7808195428
Set-Txt "R22" 6820858.896214374

# ---- Row 23 (becomes the old row24 record) ----
Set-Txt "A23" 111941129
Set-Txt "B23" 88489
Set-Txt "E23" 1962
Set-Txt "F23" "Vaddporing"
Set-Txt "G23" "Anomoporia kamtschatica"
Set-Txt "H23" "(Parmasto) Bondartseva"
Set-Txt "Q23" 466215.6937692813
Set-Txt "R23" 6820389.803268042
Set-Txt "AX23" "Bengt Oldhammer, Birgitta Kvist, Peter Turander"

# ---- Row 24 (becomes the old row25 record) ----
Set-Txt "A24" 111941043
Set-Txt "B24" 77268
Set-Txt "E24" 228912
Set-Txt "F24" "Mörk kolflarnlav"
Set-Txt "G24" "Carbonicola myrmecina"
Set-Txt "H24" "(Ach.) Bendiksby & Timdal"
Set-Txt "Q24" 466184.335225084
Set-Txt "R24" 6820409.199356439

# ---- Row 25 (becomes the old row22 record) ----
Set-Txt "A25" 111941765
Set-Txt "B25" 56398
Set-Txt "E25" 100109
Set-Txt "F25" "Tretåig hackspett"
Set-Txt "G25" "Picoides tridactylus"
Set-Txt "H25" "(Linnaeus, 1758)"
Clear-Cell "J25"
Clear-Cell "AF25"
Set-Txt "L25" ""
Set-Txt "M25" "äldre spår"
Set-Txt "Q25" 466413.7788343028
Set-Txt "R25" 6820854.133933268
Set-Txt "AX25" "Bengt Oldhammer, Peter Turander, Birgitta Kvist"

# ---- Row 29 (becomes the old row30 record) ----
Set-Txt "A29" 111941668
Set-Txt "B29" 56398
Set-Txt "E29" 100109
Set-Txt "F29" "Tretåig hackspett"
Set-Txt "G29" "Picoides tridactylus"
Set-Txt "H29" "(Linnaeus, 1758)"
Clear-Cell "J29"
Clear-Cell "AF29"
Set-Txt "L29" ""
Set-Txt "M29" "äldre spår"
Set-Txt "Q29" 466427.1314100454
Set-Txt "R29" 6820752.359779999
Set-Txt "AX29" "Bengt Oldhammer, Peter Turander, Birgitta Kvist"

# ---- Row 30 (becomes the old row29 record) ----
Set-Txt "A30" 111941321
Set-Txt "B30" 77515
Set-Txt "E30" 6425
Set-Txt "F30" "Garnlav"
Set-Txt "G30" "Alectoria sarmentosa"
Set-Txt "H30" "(Ach.) Ach."
Clear-Cell "L30"
Clear-Cell "M30"
Set-Txt "J30" ""
Set-Txt "AF30" ""
Set-Txt "Q30" 466314.5865380571
Set-Txt "R30" 6820482.860897989
Set-Txt "AX30" "Bengt Oldhammer, Birgitta Kvist, Peter Turander"
